$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: species record changes from "Tretåig hackspett" to "Blåsippa"
$ws.Range("A3").Value = 112042940
$ws.Range("B3").Value = 98961
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 222498
$ws.Range("F3").Value = "Blåsippa"
$ws.Range("G3").Value = "Hepatica nobilis"
$ws.Range("H3").Value = "Schreb."
$ws.Range("P3").Value = "Stor-Moberg (Stor-Moberg), Dlr"
$ws.Range("Q3").Value = 511611
$ws.Range("R3").Value = 6733626
$ws.Range("S3").Value = 1
$ws.Range("Z3").Value = "10:33"
$ws.Range("AB3").Value = "10:33"
$ws.Range("AC3").Value = "Fullt med blåsippsblad på denna sidan bäcken"
$ws.Range("AW3").Value = "Evalena Sköld"
$ws.Range("AX3").Value = "Evalena Sköld, Åke Sköld"

# Row 4: taxon sort order changed
$ws.Range("B4").Value = 90466

# Row 5: taxon sort order changed
$ws.Range("B5").Value = 95687

# Row 6: species record changes from "Blåsippa" to "Tretåig hackspett"
$ws.Range("A6").Value = 112044333
$ws.Range("B6").Value = 56430
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("P6").Value = "Stor Mpmerg, Kilen-Stor, Moberg, Leksand, Dlr"
$ws.Range("Q6").Value = 511614
$ws.Range("R6").Value = 6733640
$ws.Range("S6").Value = 25
$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()
$ws.Range("AC6").Value = "Minst 2"
$ws.Range("AW6").Value = "Åke Sköld"
$ws.Range("AX6").Value = "Åke Sköld"
